$d = $word.ActiveDocument

# Locate the 3rd occurrence of "chosen sample" -> this is the pathogenTest
# paragraph within the "Templates for Event Participants" section (the
# Cases and Contacts sections each contain one earlier occurrence).
$searchRng = $d.Content
$matchCount = 0
$targetRng = $null
while ($searchRng.Find.Execute("chosen sample", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $matchCount = $matchCount + 1
    if ($matchCount -eq 3) {
        $targetRng = $searchRng.Duplicate
    }
    $searchRng.Collapse(0)
    if ($matchCount -ge 10) { break }
}

if ($targetRng -eq $null) {
    throw "Could not locate the 'Event Participants' pathogenTest paragraph"
}

# The pathogenTest paragraph is followed by a blank, shaded separator
# paragraph; the new "Templates for Travel Entries" section is inserted
# right after that separator (i.e. immediately before "Templates for
# Event Handouts").
$pathogenPara = $targetRng.Paragraphs(1)
$separatorPara = $pathogenPara.Next()
$insertPos = $separatorPara.Range.End
$insertPoint = $d.Range($insertPos, $insertPos)

$newBodyXml = '<w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="9"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Templates for Travel Entries</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> (.</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>docx</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>):</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Available root entities:</w:t></w:r><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>travelEntry</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>travelEntry</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> data</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>person</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: the travelEntry person</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Listenabsatz"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="11"/></w:numPr><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF" w:themeFill="background1"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>user</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: the current user</w:t></w:r></w:p>'

# A trailing empty <w:p/> is appended because Word's InsertXML merges the
# *last* inserted paragraph into the paragraph found at the insertion
# point (taking over its paragraph mark) instead of creating a clean
# break. Appending a throwaway empty paragraph absorbs that merge so the
# four real paragraphs above it stay intact; the throwaway is then
# deleted.
$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $newBodyXml + '<w:p></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($packageXml)

# Re-find the heading we just inserted, then walk forward exactly four
# paragraphs (heading, travelEntry, person, user) to land on the
# throwaway paragraph that needs to be merged away.
$headingRng = $d.Content
$headingRng.Find.Execute("Templates for Travel Entries", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$walkPara = $headingRng.Paragraphs(1)
for ($i = 0; $i -lt 4; $i++) {
    $walkPara = $walkPara.Next()
}
$walkPara.Range.Delete()

Write-Host "Inserted Travel Entries template section at position $insertPos"
